$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "what" field of the Master's degree entry (row 2) to reflect
# that the degree is still in progress.
$ws.Range("D2").Value = "M.Sc. in Economics (in progress)"

# Leave the active selection on D3, mirroring the cursor position after
# editing D2 in the original session.
$ws.Range("D3").Select()
